$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for column C
$ws.Range("C1").Value = "Terminal Block"

# Terminal block numbers for the populated rows (mirrors the I/O list in columns A/B)
$ws.Range("C10").Value = 8
$ws.Range("C11").Value = 9
$ws.Range("C12").Value = 10

$ws.Range("C19").Value = 19
$ws.Range("C20").Value = 20
$ws.Range("C21").Value = 21
$ws.Range("C22").Value = 17
$ws.Range("C23").Value = 18
$ws.Range("C24").Value = 15
$ws.Range("C25").Value = 16
$ws.Range("C26").Value = 11
$ws.Range("C27").Value = 12
$ws.Range("C28").Value = 13
$ws.Range("C29").Value = 14
$ws.Range("C30").Value = 22

$ws.Range("C35").Value = 1
$ws.Range("C36").Value = 2
$ws.Range("C37").Value = 3
$ws.Range("C38").Value = 4
$ws.Range("C39").Value = 5
$ws.Range("C40").Value = 6
$ws.Range("C41").Value = 7

# Update the selected cell / view to match the new active cell after editing
$ws.Range("C34").Select()
